$d = $word.ActiveDocument
$pairs = @(
  ,@("[Music]", "[Muziki]")
  ,@("Introduction", "Utangulizi")
  ,@("Hi everyone I'm sony from London in the", "Jambo kila mtu mimi nina sony kutoka London katika")
  ,@("UK and today I have a really exciting", "Uingereza na leo nina kusisimua sana")
  ,@("video for you for the Virtual Maths", "video kwako kwa Hesabu pepe")
  ,@("Camp (VMC)", "Kambi (VMC)")
  ,@("We're going to play a fun mathematical game called split or steal", "Tutacheza mchezo wa kufurahisha wa hisabati unaoitwa kugawanyika au kuiba")
  ,@("and this game is a two player game.", "na mchezo huu ni wa wachezaji wawili.")
  ,@("So a quick introduction split or steal", "Kwa hivyo utangulizi wa haraka hugawanyika au kuiba")
  ,@("is a game based on the famous Prisoners Dilemma which you might have", "ni mchezo unaotokana na Dilemma maarufu ya Wafungwa ambayo unaweza kuwa nayo")
  ,@("heard of before", "kusikia kabla")
  ,@("It's an interesting new branch of maths called Game Theory", "Ni tawi jipya la kuvutia la hisabati linaloitwa Nadharia ya Mchezo")
  ,@("Where when making your choice you also need to consider", "Ambapo unapofanya uchaguzi wako pia unahitaji kuzingatia")
  ,@("The other player's choice", "Chaguo la mchezaji mwingine")
  ,@("This means your choices are interdependent ", "Hii inamaanisha kuwa chaguzi zako zinategemeana ")
  ,@("So what you need today: a partner to play the game with,", "Kwa hivyo unachohitaji leo: mshirika wa kucheza naye mchezo,")
  ,@("two small cards each saying split and steal,", "kadi mbili ndogo kila moja ikisema gawanya na uibe,")
  ,@("and these can just be small bits of paper", "na hizi zinaweza kuwa vipande vidogo vya karatasi")
  ,@("and finally something to count scores with", "na hatimaye kitu cha kuhesabu alama nacho")
  ,@("If you have some tokens or coins or something like that,", "Ikiwa una ishara au sarafu au kitu kama hicho,")
  ,@("That would be great,", "Hiyo itakuwa nzuri,")
  ,@("but if not you could just use a pen and paper to write down the", "lakini kama sivyo unaweza tu kutumia kalamu na karatasi kuandika")
  ,@("scores ", "alama ")
  ,@("So this is what the game looks and this table is called a matrix", "Hivyo hii ni nini mchezo inaonekana na meza hii inaitwa tumbo")
  ,@("In the game there are two points to be won but who wins what is decided by the choices of the players", "Katika mchezo huo kuna pointi mbili za kushinda lakini nani atashinda kile kinachoamuliwa na uchaguzi wa wachezaji")
  ,@("We have our two players on the left and at the top", "Tuna wachezaji wetu wawili kushoto na juu")
  ,@("The red player and the blue player", "Mchezaji mwekundu na mchezaji wa bluu")
  ,@("each player has two options shown next to them, split or steal", "kila mchezaji ana chaguo mbili zilizoonyeshwa karibu nao, kugawanyika au kuiba")
  ,@("Since two players have two choices each there are four outcomes in total", "Kwa kuwa wachezaji wawili wana chaguo mbili kila mmoja kuna matokeo manne kwa jumla")
  ,@("and they are all shown in the table in each section", "na zote zimeonyeshwa kwenye jedwali katika kila sehemu")
  ,@("The red number is the number of points won by the red player", "Nambari nyekundu ni idadi ya pointi alizoshinda mchezaji mwekundu")
  ,@("and the blue number is the number of points won by the blue player", "na nambari ya bluu ni nambari ya alama alizoshinda mchezaji wa bluu")
  ,@("for example", "kwa mfano")
  ,@("if both players choose to split", "ikiwa wachezaji wote wawili watachagua kugawanyika")
  ,@("we would end up with the top left outcome", "tungeishia na matokeo ya juu kushoto")
  ,@("and the players would split the two points to earn one point each", "na wachezaji wangegawanya pointi mbili ili kupata pointi moja kila mmoja")
  ,@("However if the red player wanted to split", "Walakini ikiwa mchezaji nyekundu alitaka kugawanyika")
  ,@("But the blue player chose to steal the blue player would steal the red player's point and earn two points", "Lakini mchezaji wa bluu alichagua kuiba mchezaji wa bluu angeiba pointi ya mchezaji nyekundu na kupata pointi mbili")
  ,@("while the red player wins nothing", "huku mchezaji mwekundu hashindi chochote")
  ,@("The opposite happens if the blue player splits and the red player steals", "Kinyume chake hutokea ikiwa mchezaji wa bluu atagawanyika na mchezaji nyekundu anaiba")
  ,@("but if both players try and steal", "lakini ikiwa wachezaji wote wawili watajaribu na kuiba")
  ,@("it doesn't work and no one wins the points so both players end up with nothing", "haifanyi kazi na hakuna anayeshinda pointi hivyo wachezaji wote wawili wanaishia bila chochote")
  ,@("Now it's your turn get ready to play split or steal with your partner", "Sasa ni zamu yako jitayarishe kucheza mgawanyiko au kuiba na mwenzako")
  ,@("First talk to your partner for a couple of minutes about what choice you're going to make", "Kwanza zungumza na mwenza wako kwa dakika kadhaa kuhusu chaguo utakalofanya")
  ,@("Remember, the person with the most points wins", "Kumbuka, mtu aliye na pointi nyingi hushinda")
  ,@("You are allowed to lie to your partner then secretly choose split or steal", "Unaruhusiwa kumdanganya mwenzako kisha kwa siri chagua kugawanyika au kuiba")
  ,@("and place the card you have chosen face down so your partner can't see it", "na weka kadi uliyochagua kifudifudi ili mwenzako asiione")
  ,@("Finally reveal your choices and work out your scores", "Hatimaye onyesha chaguo zako na ufanyie kazi alama zako")
  ,@("Play the game once with your partner", "Cheza mchezo mara moja na mwenzi wako")
  ,@("Pause the video now", "Sitisha video sasa")
  ,@("How did it go? ", "Iliendaje? ")
  ,@("Did you get the number of points you were hoping for?", "Je, umepata idadi ya pointi ulizokuwa unazitarajia?")
  ,@("did you and your partner tell the truth to each other?", "wewe na mwenzako mliambiana ukweli?")
  ,@("Let's think about why the result might have been different to what you expected", "Hebu tufikirie kwa nini matokeo yanaweza kuwa tofauti na ulivyotarajia")
  ,@("Imagine you are the red player your opponent", "Fikiria wewe ni mchezaji nyekundu mpinzani wako")
  ,@("The blue player has two choices split or steal ", "Mchezaji wa bluu ana chaguo mbili za kugawanyika au kuiba ")
  ,@("if blue chooses split", "ikiwa bluu inachagua kugawanyika")
  ,@("you could either choose split and win one point", "unaweza kuchagua kugawanyika na kushinda pointi moja")
  ,@("or you could choose steal and win two points", "au unaweza kuchagua kuiba na kushinda pointi mbili")
  ,@("two points is better than one so you would choose steal", "pointi mbili ni bora kuliko moja hivyo ungependa kuchagua kuiba")
  ,@("What if the blue player chose steal?", "Je, ikiwa mchezaji wa bluu atachagua kuiba?")
  ,@("If you choose split you get zero and if you choose steal you also get zero", "Ukichagua kugawanya utapata sifuri na ukichagua kuiba pia utapata sifuri")
  ,@("So it doesn't really matter what you choose", "Kwa hivyo haijalishi unachagua nini")
  ,@("But let's assume you prefer to steal so you don't give your opponent any points", "Lakini hebu tuchukulie unapendelea kuiba ili usimpe mpinzani wako pointi zozote")
  ,@("As we have seen no matter what your opponent does", "Kama tulivyoona, haijalishi mpinzani wako anafanya nini")
  ,@("split is never the best choice this means steal is called a weakly dominant strategy", "mgawanyiko kamwe sio chaguo bora zaidi hii inamaanisha kuiba inaitwa mkakati dhaifu")
  ,@("Since this is a symmetric game steel is also weakly dominant for the blue player", "Kwa kuwa huu ni mchezo wa chuma wa ulinganifu pia hutawala kwa njia dhaifu kwa kicheza bluu")
  ,@("We have proven that it makes sense for both players to choose steal", "Tumethibitisha kwamba inaleta maana kwa wachezaji wote wawili kuchagua kuiba")
  ,@("Therefore the steel steel outcome is known as the nash equilibrium", "Kwa hivyo matokeo ya chuma ya chuma hujulikana kama usawa wa nash")
  ,@("but look at the matrix the split split outcome is better for both players as they both get one point instead of zero ", "lakini angalia matrix matokeo ya mgawanyiko wa mgawanyiko ni bora kwa wachezaji wote wawili kwani wote wanapata alama moja badala ya sifuri ")
  ,@("This means the nash equilibrium is not the optimal solution, an incredible result", "Hii inamaanisha kuwa usawa wa nash sio suluhisho bora, matokeo ya kushangaza")
  ,@("We now know what should happen for a single game", "Sasa tunajua nini kifanyike kwa mchezo mmoja")
  ,@("But does this result hold if we play the game multiple times", "Lakini je, matokeo haya yanashikilia ikiwa tutacheza mchezo mara nyingi")
  ,@("against the same player", "dhidi ya mchezaji huyo huyo")
  ,@("Start the scoring from zero and play the game ten times in a row", "Anza kufunga kutoka sifuri na ucheze mchezo mara kumi mfululizo")
  ,@("and play the game ten times in a row with your partner", "na cheza mchezo huo mara kumi mfululizo na mwenzako")
  ,@("Does your strategy change?", "Je, mkakati wako unabadilika?")
  ,@("Now you know", "Sasa unajua")
  ,@("You'll be playing the same opponent again", "Utakuwa unacheza mpinzani sawa tena")
  ,@("Play the game 10 times with the same partner", "Cheza mchezo mara 10 na mwenzi sawa")
  ,@("[PAUSE]", "[SItisha]")
  ,@("Did you manage to score more points than your opponent?", "Je, umeweza kupata pointi zaidi ya mpinzani wako?")
  ,@("A repeated game like the one you've just played is much more complicated", "Mchezo unaorudiwa kama ule ambao umecheza hivi punde ni mgumu zaidi")
  ,@("Because your decision is not only influenced by your communication with your partner", "Kwa sababu uamuzi wako hauathiriwi tu na mawasiliano yako na mwenza wako")
  ,@("In this round but also what has happened in previous rounds ", "Katika raundi hii lakini pia yale yaliyotokea katika raundi zilizopita ")
  ,@("You might trust your partner less if they stole in the previous round which could make you ", "Unaweza kumwamini mwenzi wako kidogo ikiwa aliiba katika raundi ya awali ambayo inaweza kukufanya ")
  ,@("more likely to steal in this round", "uwezekano mkubwa wa kuiba katika raundi hii")
  ,@("In general,", "Kwa ujumla,")
  ,@("The more the game is repeated the more likely you will be to cooperate with your opponent", "Kadiri mchezo unavyorudiwa ndivyo unavyoweza kushirikiana na mpinzani wako")
  ,@("Because they could punish you in future rounds if you don't", "Kwa sababu wanaweza kukuadhibu katika raundi zijazo usipofanya hivyo")
  ,@("Using the same logic, if you know there aren't many rounds left", "Kwa kutumia mantiki sawa, ikiwa unajua hakuna raundi nyingi zilizobaki")
  ,@("You might be more tempted to steal because your opponent has less time to retaliate", "Huenda ukajaribiwa zaidi kuiba kwa sababu mpinzani wako ana muda mchache wa kulipiza kisasi")
  ,@("As this is a very famous game game theorists have developed many strategies", "Kama huu ni mchezo maarufu sana wananadharia wameunda mikakati mingi")
  ,@("that we could use when playing ", "ambayo tunaweza kutumia wakati wa kucheza ")
  ,@("You could always cooperate, meaning choosing split every time", "Unaweza kushirikiana kila wakati, ikimaanisha kuchagua mgawanyiko kila wakati")
  ,@("or you could play steal every time", "au unaweza kucheza kuiba kila wakati")
  ,@("You might choose to copy what your opponent did in their last move", "Unaweza kuchagua kunakili kile mpinzani wako alifanya katika hatua yao ya mwisho")
  ,@("sometimes known as tit for tat or copycat", "wakati mwingine hujulikana kama tit kwa tat au copycat")
  ,@("Grim trigger is where you play split but if your opponent plays steal just once you punish them by playing steal for the rest of the game", "Grim trigger ni pale unapocheza mgawanyiko lakini mpinzani wako akicheza kuiba mara moja tu, unamwadhibu kwa kucheza kuiba kwa muda wote uliosalia wa mchezo")
  ,@("You could even decide to choose randomly each time by flipping a coin", "Unaweza hata kuamua kuchagua nasibu kila wakati kwa kugeuza sarafu")
  ,@("Which strategy do you think is best? ", "Je, unadhani ni mkakati gani bora zaidi? ")
  ,@("Try playing five rounds sticking to one of the strategies listed", "Jaribu kucheza raundi tano ukitumia moja ya mikakati iliyoorodheshwa")
  ,@("and see what happens if you can swap partners this time", "na uone kitakachotokea ikiwa unaweza kubadilishana washirika wakati huu")
  ,@("Then pick another strategy and play five more rounds", "Kisha chagua mkakati mwingine na ucheze raundi nyingine tano")
  ,@("play using a strategy ", "kucheza kwa kutumia mkakati ")
  ,@("So, which strategy scored you the most points", "Kwa hivyo, ni mkakati gani ulikupa alama nyingi zaidi")
  ,@("In 1980, Robert Axelrod made a tournament where he played 63 different strategies against each other", "Mnamo 1980, Robert Axelrod alifanya mashindano ambapo alicheza mikakati 63 tofauti dhidi ya kila mmoja")
  ,@("To see which one came out on top", "Ili kuona ni yupi aliyetoka juu")
  ,@("and out of all of them it was tit for tat that won", "na kati ya hao wote ilikuwa tit kwa tat iliyoshinda")
  ,@("In general the most successful strategies were nice, ", "Kwa ujumla mikakati iliyofanikiwa zaidi ilikuwa nzuri, ")
  ,@("Meaning they started off cooperating by playing split and forgiving ", "Ikimaanisha walianza kushirikiana kwa kucheza kugawanyika na kusameheana ")
  ,@("Meaning that they wouldn't do what grim trigger does", "Ikimaanisha kuwa hawangefanya kile kichochezi kibaya hufanya")
  ,@("and fully stop cooperating once the opponent played steal", "na kuacha kabisa kushirikiana mara tu mpinzani alipocheza kuiba")
  ,@("I guess the fact that nice and forgiving strategies are the best", "Nadhani ukweli kwamba mikakati mizuri na ya kusamehe ni bora zaidi")
  ,@("is a good sign for society", "ni ishara nzuri kwa jamii")
  ,@("and that's the end of this session", "na huo ndio mwisho wa kikao hiki")
  ,@("if you enjoyed the topic, there's a very good website", "ikiwa ulifurahia mada, kuna tovuti nzuri sana")
  ,@("called nikki case's evolution of trust", "inayoitwa nikki case's evolution of trust")
  ,@("which goes into more detail", "ambayo inaingia kwa undani zaidi")
  ,@("Thank you and enjoy the rest of your Virtual Maths Camp.", "Asante na ufurahie mapumziko ya Kambi yako ya Hisabati Mtandaoni.")
)

$totalReplacements = 0
$missed = @()
foreach ($pair in $pairs) {
  $old = $pair[0]
  $new = $pair[1]
  $r = $d.Content
  $countThis = 0
  while ($r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r.Text = $new
    $countThis = $countThis + 1
    $totalReplacements = $totalReplacements + 1
    $r.Collapse(0)
    $r.End = $d.Content.End
  }
  if ($countThis -eq 0) {
    $missed += $old
  }
}

Write-Output "Total replacements: $totalReplacements"
Write-Output "Missed count: $($missed.Count)"
foreach ($m in $missed) { Write-Output "MISSED: $m" }
